$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (H1, style index 1: bold / border / centered)
# onto the two new header cells so I1/J1 match the existing header styling.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iVals = @(9,9,7,6,9,8,9,9,7,10,8,6,7,7,9,7,7,8,7,6,7,6,10,7,7,6,6,9,7,4,6,8,5,8,7,7,10,6,3,7,8,6,5,8,6,6,9,7,1,4,5,6,5,9,5,5,6,8,8,2,5,4,6,3,3,6)
$jVals = @(9,9,7,6,9,8,10,9,7,10,8,7,7,7,9,7,7,8,7,6,7,7,10,8,7,6,6,10,8,4,7,8,5,8,7,8,10,6,3,7,8,7,5,9,6,6,9,7,2,4,7,6,5,9,5,6,7,8,8,2,5,4,6,3,3,6)

for ($r = 2; $r -le 67; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
